# Updated cryptos list on Tue May 21 19:58:19 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.377.89'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.53%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.714.83'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +8.40%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '609.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.13'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.43%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.711.61'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +7.92%  '
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("E9").Value = '  +0.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.165'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.62%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.35'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.496'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.47%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.47'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.16%  '
$ws.Range("E14").Value = '  +0.97%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.334.76'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +8.56%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.713.12'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +8.90%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.513.05'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.24%  '
$ws.Range("E18").Value = '  +0.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.55'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.86%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '512.13'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.35%  '
$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.68'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.39%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.42'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +12.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.724'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.47%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '87.41'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.46'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.89%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.31'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.93'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.18%  '
$ws.Range("E29").Value = '  +19.90%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.50'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.64%  '
$ws.Range("E31").Value = '  -5.18%  '
$ws.Range("E32").Value = '  +4.49%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.84%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.115'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.91%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.15'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.18%  '
$ws.Range("E37").Value = '  +2.31%  '
$ws.Range("E38").Value = '  +1.65%  '
$ws.Range("E39").Value = '  +2.90%  '
$ws.Range("E40").Value = '  +1.85%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '51.15'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '44.15'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -10.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.76'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.85%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.080.54'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '418.37'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.69%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.70'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.53%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0362'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.68'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.81%  '
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.49'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.24%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '134.55'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.73%  '
